# Amun's Book review doc:
#  1) Insert a new "Meta description" paragraph right after the title
#     (Heading1) paragraph.
#  2) Remove the duplicate bold title paragraph that was near the end of
#     the document (its content now lives in the new Meta description
#     paragraph instead).
#  3) Replace the text of the remaining (italic) paragraph at the end
#     with the new image-generation prompt text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: insert the "Meta description" paragraph right after the title.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start

# Type the whole sentence first ...
$insertRange = $d.Range($metaStart, $metaStart)
$insertRange.InsertAfter("Meta description: Experience the mystery and excitement of ancient Egypt with Amun's Book. Play this slot game for free and trigger the exciting free spins feature.")

# ... then bold just the "Meta description" label.
$labelEnd = $metaStart + ("Meta description").Length
$labelRange = $d.Range($metaStart, $labelEnd)
$labelRange.Font.Bold = 1

# Leave a leading (empty) run marker in front of the label, matching the
# style used throughout the rest of the document's body paragraphs.
$leadRange = $d.Range($metaStart, $metaStart)
$leadRange.InsertBefore("")

# ---------------------------------------------------------------------
# Step 2: delete the old duplicate bold title paragraph near the end of
# the document (second-to-last paragraph at this point).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
if ($dupTitlePara.Range.Text.Trim() -eq "Play Amun's Book Free - Exciting Egyptian-themed Slot Game") {
    $dupTitlePara.Range.Delete()
} else {
    Write-Host "WARNING: unexpected paragraph, skipping delete: [$($dupTitlePara.Range.Text)]"
}

# ---------------------------------------------------------------------
# Step 3: update the final (italic) paragraph's text to the new prompt.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($count)
$pStart = $promptPara.Range.Start
$pEnd = $promptPara.Range.End

# Exclude the trailing paragraph-mark character from the replacement.
$promptContent = $d.Range($pStart, $pEnd - 1)
$promptContent.Text = "Prompt: Create a cartoon-style feature image that will grab the attention of slot players for the game `"Amun's Book`". The image should feature a happy Maya warrior with glasses. The Maya warrior should be happy and holding up a book with the title `"Amun's Book`" written on it. The background of the image should be a dark Egyptian temple with stone columns and hieroglyphics in the background. The temple should be lit with torches to create a mysterious and adventurous feel. Additionally, include some of the slot's elements, such as the scarab symbol or the book of Amun, to give players a sense of what the game is about. The overall image should have bright and vibrant colors to make it stand out and create a fun atmosphere for players."

Write-Host "Meta description paragraph: [$($metaPara.Range.Text)]"
Write-Host "Prompt paragraph: [$($promptPara.Range.Text)]"
Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
